$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.088.88"
$ws.Range("E2").Value = "  +1.01%  "

$ws.Range("D3").Value = "2.345.10"
$ws.Range("E3").Value = "  +5.14%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.34"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.13"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.76%  "

$ws.Range("E7").Value = "  +1.75%  "

$ws.Range("E8").Value = "  -0.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.622"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.36"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.97"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.07"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +19.29%  "

$ws.Range("E14").Value = "  +0.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.19"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +6.97%  "

$ws.Range("D16").Value = "2.704.99"
$ws.Range("E16").Value = "  +5.00%  "

$ws.Range("D17").Value = "2.424.91"
$ws.Range("E17").Value = "  +7.06%  "

$ws.Range("D18").Value = "43.081.56"
$ws.Range("E18").Value = "  +1.17%  "

$ws.Range("E19").Value = "  +1.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.24"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.57"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.45"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.36%  "

$ws.Range("E23").Value = "  +9.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "250.31"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +8.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.96"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.87%  "

$ws.Range("E26").Value = "  +0.62%  "

$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.24"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.55"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.44"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.88"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.48%  "

$ws.Range("E32").Value = "  -2.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0907"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.78"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.95"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.97%  "

$ws.Range("E36").Value = "  +3.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0377"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.68%  "

$ws.Range("E38").Value = "  -2.96%  "

$ws.Range("E39").Value = "  +0.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.75"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +9.40%  "

$ws.Range("E41").Value = "  +13.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.12"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.233"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.69%  "

$ws.Range("E44").Value = "  -0.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.55"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.69"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.19"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.43"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.93%  "

$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0998"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.54"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.10%  "

